$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 27 ("Vega Monumental Concepción"
# weekly price update). This shifts the existing rows 27-29 down to 28-30.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with this week's data.
$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 44644
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100112037
$ws.Range("G27").Value = "Cebollín"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 160
$ws.Range("K27").Value = 6500
$ws.Range("L27").Value = 7000
$ws.Range("M27").Value = 6750
$ws.Range("N27").Value = "`$/paquete 36 unidades"
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 188
$ws.Range("Q27").Value = 36
$ws.Range("R27").Value = "Hortaliza"
